# Natmi following Dr Hou advice
# Update the Ceacam1-Sele NATMI LR-pairs sheet to include the new "FAPs"
# cluster as both a sending and a target cluster, with refreshed stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the 20 column values (A:T) for one data row.
# Row order below is written starting at worksheet row 2.
$data = @(
    @("ECs", "Ceacam1", "Sele", "ECs", 3, 1, 42.924731, 128.774193, 0.8780494305754949, 0.878049430575495, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 201.7250594993734, 1815.52553549436, 0.8482392350028691, 0.8482392350028691),
    @("ECs", "Ceacam1", "Sele", "FAPs", 3, 1, 42.924731, 128.774193, 0.8780494305754949, 0.878049430575495, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 7.089348414254334, 63.804135728289, 0.02981019557262582, 0.02981019557262582),
    @("FAPs", "Ceacam1", "Sele", "ECs", 2, 0.6666666666666666, 0.1486296666666667, 0.445889, 0.003040303133950735, 0.003040303133950736, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 0.6984861093644444, 6.28637498428, 0.002937083397262636, 0.002937083397262636),
    @("FAPs", "Ceacam1", "Sele", "FAPs", 2, 0.6666666666666666, 0.1486296666666667, 0.445889, 0.003040303133950735, 0.003040303133950736, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 0.02454732894411111, 0.220925960497, 0.0001032197366880999, 0.0001032197366880999),
    @("M2", "Ceacam1", "Sele", "ECs", 3, 1, 5.009066, 15.027198, 0.1024632524549792, 0.1024632524549792, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 23.54013906077333, 211.86125154696, 0.09898457632544934, 0.09898457632544934),
    @("M2", "Ceacam1", "Sele", "FAPs", 3, 1, 5.009066, 15.027198, 0.1024632524549792, 0.1024632524549792, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 0.8272856527393333, 7.445570874653999, 0.003478676129529862, 0.003478676129529863),
    @("sCs", "Ceacam1", "Sele", "ECs", 3, 1, 0.8040363333333334, 2.412109, 0.01644701383557516, 0.01644701383557517, 3, 1, 4.699506666666667, 14.09852, 0.9660495246229048, 0.9660495246229047, 3.778574108742223, 34.00716697868, 0.01588862989732373, 0.01588862989732373),
    @("sCs", "Ceacam1", "Sele", "FAPs", 3, 1, 0.8040363333333334, 2.412109, 0.01644701383557516, 0.01644701383557517, 1, 0.3333333333333333, 0.1651576666666667, 0.495473, 0.03395047537709522, 0.03395047537709522, 0.1327927647285556, 1.195134882557, 0.000558383938251439, 0.0005583839382514391)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $targetRow = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $rowValues[$j]
    }
}
